$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new customer record (phone 79174441, blank birthday, 0 points)
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "79174441"
$ws.Cells.Item(5, 1).Style = "Normal"

$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 2).Style = "Normal"

$ws.Cells.Item(5, 3).Value = 0
